$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = "1x$($i)=$($i)"
}
